# Wrap the 'dob' value of each transaction-record string in column C
# (rows 2-128) with a single-element list literal, e.g.
#   'dob': 1992-11-10 00:00:00}   ->  'dob': ['1992-11-10 00:00:00']}
#   'dob': n/a}                   ->  'dob': ['n/a']}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $val = $cell.Text

    if ($val -match "'dob': (.*)\}\s*$") {
        $dobValue = $Matches[1]
        # Only wrap if not already wrapped in a list
        if ($dobValue -notmatch "^\[.*\]$") {
            $newVal = $val -replace "'dob': (.*)\}\s*$", "'dob': ['`$1']}"
            $cell.Value = $newVal
        }
    }
}
